$wb = $excel.ActiveWorkbook

# Insert a new "TaxonRelation" worksheet right after "Synonym" (and before "Distribution")
$synonym = $wb.Worksheets.Item("Synonym")
$taxonRelation = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $synonym)
$taxonRelation.Name = "TaxonRelation"

# Header row, mirroring the layout of the other *Relation sheets (e.g. NameRelation)
$taxonRelation.Range("A1").Value = "taxonID"
$taxonRelation.Range("B1").Value = "relatedTaxonID"
$taxonRelation.Range("C1").Value = "type"
$taxonRelation.Range("D1").Value = "referenceID"
$taxonRelation.Range("E1").Value = "remarks"

# Match the look & feel of the neighbouring sheets
$taxonRelation.Range("A1:E1").EntireColumn.ColumnWidth = 46.1640625

$taxonRelation.Activate()
$excel.ActiveWindow.Zoom = 160
